# Update NATMI LR-pair TPM-derived values (Sema3f-Plxna3), rows 2-17
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"40.83537666666667"
$ws.Range("H2").Value = [double]"122.50613"
$ws.Range("I2").Value = [double]"0.9274830900091532"
$ws.Range("J2").Value = [double]"0.9274830900091531"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.4286376666666666"
$ws.Range("N2").Value = [double]"1.285913"
$ws.Range("O2").Value = [double]"0.031593034723905"
$ws.Range("P2").Value = [double]"0.03159303472390499"
$ws.Range("Q2").Value = [double]"17.50358057185445"
$ws.Range("R2").Value = [double]"157.53222514669"
$ws.Range("S2").Value = [double]"0.02930200546849388"
$ws.Range("T2").Value = [double]"0.02930200546849387"
$ws.Range("G3").Value = [double]"40.83537666666667"
$ws.Range("H3").Value = [double]"122.50613"
$ws.Range("I3").Value = [double]"0.9274830900091532"
$ws.Range("J3").Value = [double]"0.9274830900091531"
$ws.Range("O3").Value = [double]"0.4539624805119086"
$ws.Range("P3").Value = [double]"0.4539624805119086"
$ws.Range("Q3").Value = [double]"251.5101484767067"
$ws.Range("R3").Value = [double]"2263.59133629036"
$ws.Range("S3").Value = [double]"0.4210425241734049"
$ws.Range("T3").Value = [double]"0.4210425241734049"
$ws.Range("G4").Value = [double]"40.83537666666667"
$ws.Range("H4").Value = [double]"122.50613"
$ws.Range("I4").Value = [double]"0.9274830900091532"
$ws.Range("J4").Value = [double]"0.9274830900091531"
$ws.Range("M4").Value = [double]"5.427178000000001"
$ws.Range("N4").Value = [double]"16.281534"
$ws.Range("O4").Value = [double]"0.400013895979308"
$ws.Range("P4").Value = [double]"0.400013895979308"
$ws.Range("Q4").Value = [double]"221.6208578670467"
$ws.Range("R4").Value = [double]"1994.58772080342"
$ws.Range("S4").Value = [double]"0.3710061242894886"
$ws.Range("T4").Value = [double]"0.3710061242894885"
$ws.Range("G5").Value = [double]"40.83537666666667"
$ws.Range("H5").Value = [double]"122.50613"
$ws.Range("I5").Value = [double]"0.9274830900091532"
$ws.Range("J5").Value = [double]"0.9274830900091531"
$ws.Range("M5").Value = [double]"1.552534"
$ws.Range("N5").Value = [double]"4.657602"
$ws.Range("O5").Value = [double]"0.1144305887848784"
$ws.Range("P5").Value = [double]"0.1144305887848784"
$ws.Range("Q5").Value = [double]"63.39831067780666"
$ws.Range("R5").Value = [double]"570.5847961002601"
$ws.Range("S5").Value = [double]"0.1061324360777658"
$ws.Range("T5").Value = [double]"0.1061324360777658"
$ws.Range("I6").Value = [double]"0.03813623414934058"
$ws.Range("J6").Value = [double]"0.03813623414934057"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"0.4286376666666666"
$ws.Range("N6").Value = [double]"1.285913"
$ws.Range("O6").Value = [double]"0.031593034723905"
$ws.Range("P6").Value = [double]"0.03159303472390499"
$ws.Range("Q6").Value = [double]"0.7197119325738889"
$ws.Range("R6").Value = [double]"6.477407393165"
$ws.Range("S6").Value = [double]"0.001204839369719089"
$ws.Range("T6").Value = [double]"0.001204839369719088"
$ws.Range("I7").Value = [double]"0.03813623414934058"
$ws.Range("J7").Value = [double]"0.03813623414934057"
$ws.Range("O7").Value = [double]"0.4539624805119086"
$ws.Range("P7").Value = [double]"0.4539624805119086"
$ws.Range("S7").Value = [double]"0.01731241945181761"
$ws.Range("T7").Value = [double]"0.0173124194518176"
$ws.Range("I8").Value = [double]"0.03813623414934058"
$ws.Range("J8").Value = [double]"0.03813623414934057"
$ws.Range("M8").Value = [double]"5.427178000000001"
$ws.Range("N8").Value = [double]"16.281534"
$ws.Range("O8").Value = [double]"0.400013895979308"
$ws.Range("P8").Value = [double]"0.400013895979308"
$ws.Range("Q8").Value = [double]"9.112602719163334"
$ws.Range("R8").Value = [double]"82.01342447247001"
$ws.Range("S8").Value = [double]"0.01525502360005686"
$ws.Range("T8").Value = [double]"0.01525502360005685"
$ws.Range("I9").Value = [double]"0.03813623414934058"
$ws.Range("J9").Value = [double]"0.03813623414934057"
$ws.Range("M9").Value = [double]"1.552534"
$ws.Range("N9").Value = [double]"4.657602"
$ws.Range("O9").Value = [double]"0.1144305887848784"
$ws.Range("P9").Value = [double]"0.1144305887848784"
$ws.Range("Q9").Value = [double]"2.606810675823333"
$ws.Range("R9").Value = [double]"23.46129608241"
$ws.Range("S9").Value = [double]"0.004363951727747029"
$ws.Range("T9").Value = [double]"0.004363951727747028"
$ws.Range("G10").Value = [double]"1.503819"
$ws.Range("H10").Value = [double]"4.511457"
$ws.Range("I10").Value = [double]"0.03415584247746153"
$ws.Range("J10").Value = [double]"0.03415584247746152"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"0.4286376666666666"
$ws.Range("N10").Value = [double]"1.285913"
$ws.Range("O10").Value = [double]"0.031593034723905"
$ws.Range("P10").Value = [double]"0.03159303472390499"
$ws.Range("Q10").Value = [double]"0.644593467249"
$ws.Range("R10").Value = [double]"5.801341205240999"
$ws.Range("S10").Value = [double]"0.001079086717414671"
$ws.Range("T10").Value = [double]"0.001079086717414671"
$ws.Range("G11").Value = [double]"1.503819"
$ws.Range("H11").Value = [double]"4.511457"
$ws.Range("I11").Value = [double]"0.03415584247746153"
$ws.Range("J11").Value = [double]"0.03415584247746152"
$ws.Range("O11").Value = [double]"0.4539624805119086"
$ws.Range("P11").Value = [double]"0.4539624805119086"
$ws.Range("Q11").Value = [double]"9.262207694555999"
$ws.Range("R11").Value = [double]"83.359869251004"
$ws.Range("S11").Value = [double]"0.01550547097504245"
$ws.Range("T11").Value = [double]"0.01550547097504244"
$ws.Range("G12").Value = [double]"1.503819"
$ws.Range("H12").Value = [double]"4.511457"
$ws.Range("I12").Value = [double]"0.03415584247746153"
$ws.Range("J12").Value = [double]"0.03415584247746152"
$ws.Range("M12").Value = [double]"5.427178000000001"
$ws.Range("N12").Value = [double]"16.281534"
$ws.Range("O12").Value = [double]"0.400013895979308"
$ws.Range("P12").Value = [double]"0.400013895979308"
$ws.Range("Q12").Value = [double]"8.161493392782001"
$ws.Range("R12").Value = [double]"73.453440535038"
$ws.Range("S12").Value = [double]"0.01366281161986493"
$ws.Range("T12").Value = [double]"0.01366281161986492"
$ws.Range("G13").Value = [double]"1.503819"
$ws.Range("H13").Value = [double]"4.511457"
$ws.Range("I13").Value = [double]"0.03415584247746153"
$ws.Range("J13").Value = [double]"0.03415584247746152"
$ws.Range("M13").Value = [double]"1.552534"
$ws.Range("N13").Value = [double]"4.657602"
$ws.Range("O13").Value = [double]"0.1144305887848784"
$ws.Range("P13").Value = [double]"0.1144305887848784"
$ws.Range("Q13").Value = [double]"2.334730127346"
$ws.Range("R13").Value = [double]"21.012571146114"
$ws.Range("S13").Value = [double]"0.003908473165139483"
$ws.Range("T13").Value = [double]"0.003908473165139482"
$ws.Range("E14").Value = [double]"1"
$ws.Range("F14").Value = [double]"0.3333333333333333"
$ws.Range("G14").Value = [double]"0.009899"
$ws.Range("H14").Value = [double]"0.029697"
$ws.Range("I14").Value = [double]"0.0002248333640447365"
$ws.Range("J14").Value = [double]"0.0002248333640447365"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"0.4286376666666666"
$ws.Range("N14").Value = [double]"1.285913"
$ws.Range("O14").Value = [double]"0.031593034723905"
$ws.Range("P14").Value = [double]"0.03159303472390499"
$ws.Range("Q14").Value = [double]"0.004243084262333333"
$ws.Range("R14").Value = [double]"0.038187758361"
$ws.Range("S14").Value = [double]"7.103168277357735E-06"
$ws.Range("T14").Value = [double]"7.103168277357733E-06"
$ws.Range("E15").Value = [double]"1"
$ws.Range("F15").Value = [double]"0.3333333333333333"
$ws.Range("G15").Value = [double]"0.009899"
$ws.Range("H15").Value = [double]"0.029697"
$ws.Range("I15").Value = [double]"0.0002248333640447365"
$ws.Range("J15").Value = [double]"0.0002248333640447365"
$ws.Range("O15").Value = [double]"0.4539624805119086"
$ws.Range("P15").Value = [double]"0.4539624805119086"
$ws.Range("Q15").Value = [double]"0.06096916847599999"
$ws.Range("R15").Value = [double]"0.548722516284"
$ws.Range("S15").Value = [double]"0.0001020659116435856"
$ws.Range("T15").Value = [double]"0.0001020659116435856"
$ws.Range("E16").Value = [double]"1"
$ws.Range("F16").Value = [double]"0.3333333333333333"
$ws.Range("G16").Value = [double]"0.009899"
$ws.Range("H16").Value = [double]"0.029697"
$ws.Range("I16").Value = [double]"0.0002248333640447365"
$ws.Range("J16").Value = [double]"0.0002248333640447365"
$ws.Range("M16").Value = [double]"5.427178000000001"
$ws.Range("N16").Value = [double]"16.281534"
$ws.Range("O16").Value = [double]"0.400013895979308"
$ws.Range("P16").Value = [double]"0.400013895979308"
$ws.Range("Q16").Value = [double]"0.053723635022"
$ws.Range("R16").Value = [double]"0.483512715198"
$ws.Range("S16").Value = [double]"8.993646989766914E-05"
$ws.Range("T16").Value = [double]"8.993646989766912E-05"
$ws.Range("E17").Value = [double]"1"
$ws.Range("F17").Value = [double]"0.3333333333333333"
$ws.Range("G17").Value = [double]"0.009899"
$ws.Range("H17").Value = [double]"0.029697"
$ws.Range("I17").Value = [double]"0.0002248333640447365"
$ws.Range("J17").Value = [double]"0.0002248333640447365"
$ws.Range("M17").Value = [double]"1.552534"
$ws.Range("N17").Value = [double]"4.657602"
$ws.Range("O17").Value = [double]"0.1144305887848784"
$ws.Range("P17").Value = [double]"0.1144305887848784"
$ws.Range("Q17").Value = [double]"0.015368534066"
$ws.Range("R17").Value = [double]"0.138316806594"
$ws.Range("S17").Value = [double]"2.572781422612412E-05"
$ws.Range("T17").Value = [double]"2.572781422612411E-05"
